$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete row 6 (data now ends at row 5)
$ws.Rows.Item(6).Delete()

# 2. Update column widths (OOXML width target values achieved via the
#    ColumnWidth -> width quantization of this engine: using X + 1/7
#    lands squarely in the window that rounds to exactly (X+1).
$ws.Columns.Item(2).ColumnWidth = 7.142857142857143
$ws.Columns.Item(3).ColumnWidth = 7.142857142857143
$ws.Columns.Item(5).ColumnWidth = 7.142857142857143
$ws.Columns.Item(6).ColumnWidth = 7.142857142857143
$ws.Columns.Item(7).ColumnWidth = 7.142857142857143
$ws.Columns.Item(9).ColumnWidth = 7.142857142857143
$ws.Columns.Item(10).ColumnWidth = 7.142857142857143
$ws.Columns.Item(11).ColumnWidth = 7.142857142857143
$ws.Columns.Item(12).ColumnWidth = 7.142857142857143
$ws.Columns.Item(13).ColumnWidth = 7.142857142857143
$ws.Columns.Item(15).ColumnWidth = 7.142857142857143
$ws.Columns.Item(16).ColumnWidth = 7.142857142857143
$ws.Columns.Item(17).ColumnWidth = 7.142857142857143
$ws.Columns.Item(20).ColumnWidth = 8.142857142857142
$ws.Columns.Item(21).ColumnWidth = 7.142857142857143
$ws.Columns.Item(22).ColumnWidth = 7.142857142857143
$ws.Columns.Item(23).ColumnWidth = 7.142857142857143
$ws.Columns.Item(24).ColumnWidth = 7.142857142857143
$ws.Columns.Item(26).ColumnWidth = 7.142857142857143
$ws.Columns.Item(27).ColumnWidth = 7.142857142857143
$ws.Columns.Item(28).ColumnWidth = 7.142857142857143
$ws.Columns.Item(29).ColumnWidth = 7.142857142857143
$ws.Columns.Item(30).ColumnWidth = 7.142857142857143
$ws.Columns.Item(34).ColumnWidth = 7.142857142857143

# 3. Overwrite data rows 2-5 with the new dataset values
# Row 2
$ws.Cells.Item(2, 1).Value = 45099.50694444445
$ws.Cells.Item(2, 2).Value = 14.835
$ws.Cells.Item(2, 3).Value = 9.791
$ws.Cells.Item(2, 4).Value = 3.698
$ws.Cells.Item(2, 5).Value = 32.243
$ws.Cells.Item(2, 6).Value = 24.166
$ws.Cells.Item(2, 7).Value = 11.51
$ws.Cells.Item(2, 8).Value = 34.958
$ws.Cells.Item(2, 9).Value = 18.033
$ws.Cells.Item(2, 10).Value = 7.29
$ws.Cells.Item(2, 11).Value = 10.735
$ws.Cells.Item(2, 12).Value = 12.533
$ws.Cells.Item(2, 13).Value = 13.25
$ws.Cells.Item(2, 14).Value = 3.739
$ws.Cells.Item(2, 15).Value = 11.655
$ws.Cells.Item(2, 16).Value = 16.06
$ws.Cells.Item(2, 17).Value = 10.282
$ws.Cells.Item(2, 18).Value = 3.096
$ws.Cells.Item(2, 19).Value = 1.74
$ws.Cells.Item(2, 20).Value = 170.025
$ws.Cells.Item(2, 21).Value = 32.298
$ws.Cells.Item(2, 22).Value = 10.758
$ws.Cells.Item(2, 23).Value = 20.812
$ws.Cells.Item(2, 24).Value = 10.713
$ws.Cells.Item(2, 25).Value = 2.837
$ws.Cells.Item(2, 26).Value = 18.288
$ws.Cells.Item(2, 27).Value = 9.502000000000001
$ws.Cells.Item(2, 28).Value = 8.641999999999999
$ws.Cells.Item(2, 29).Value = 10.303
$ws.Cells.Item(2, 30).Value = 12.679
$ws.Cells.Item(2, 31).Value = 3.311
$ws.Cells.Item(2, 32).Value = 31.418
$ws.Cells.Item(2, 33).Value = 5.68
$ws.Cells.Item(2, 34).Value = 13.449
# Row 3
$ws.Cells.Item(3, 1).Value = 45099.51388888889
$ws.Cells.Item(3, 2).Value = 19.164
$ws.Cells.Item(3, 3).Value = 13.795
$ws.Cells.Item(3, 4).Value = 1.87
$ws.Cells.Item(3, 5).Value = 41.971
$ws.Cells.Item(3, 6).Value = 33.331
$ws.Cells.Item(3, 7).Value = 14.963
$ws.Cells.Item(3, 8).Value = 56.57
$ws.Cells.Item(3, 9).Value = 23.269
$ws.Cells.Item(3, 10).Value = 10.102
$ws.Cells.Item(3, 11).Value = 14.811
$ws.Cells.Item(3, 12).Value = 16.667
$ws.Cells.Item(3, 13).Value = 17.669
$ws.Cells.Item(3, 14).Value = 4.831
$ws.Cells.Item(3, 15).Value = 15.038
$ws.Cells.Item(3, 16).Value = 21.216
$ws.Cells.Item(3, 17).Value = 12.945
$ws.Cells.Item(3, 18).Value = 1.488
$ws.Cells.Item(3, 19).Value = 1.106
$ws.Cells.Item(3, 20).Value = 221.583
$ws.Cells.Item(3, 21).Value = 42.043
$ws.Cells.Item(3, 22).Value = 13.881
$ws.Cells.Item(3, 23).Value = 27.897
$ws.Cells.Item(3, 24).Value = 14.583
$ws.Cells.Item(3, 25).Value = 2.597
$ws.Cells.Item(3, 26).Value = 27.98
$ws.Cells.Item(3, 27).Value = 12.261
$ws.Cells.Item(3, 28).Value = 10.997
$ws.Cells.Item(3, 29).Value = 12.968
$ws.Cells.Item(3, 30).Value = 17.232
$ws.Cells.Item(3, 31).Value = 1.246
$ws.Cells.Item(3, 32).Value = 51.768
$ws.Cells.Item(3, 33).Value = 7.647
$ws.Cells.Item(3, 34).Value = 17.354
# Row 4
$ws.Cells.Item(4, 1).Value = 45099.52083333334
$ws.Cells.Item(4, 2).Value = 8.598000000000001
$ws.Cells.Item(4, 3).Value = 6.09
$ws.Cells.Item(4, 4).Value = 1.048
$ws.Cells.Item(4, 5).Value = 18.963
$ws.Cells.Item(4, 6).Value = 14.726
$ws.Cells.Item(4, 7).Value = 6.677
$ws.Cells.Item(4, 8).Value = 29.839
$ws.Cells.Item(4, 9).Value = 10.471
$ws.Cells.Item(4, 10).Value = 4.499
$ws.Cells.Item(4, 11).Value = 6.456
$ws.Cells.Item(4, 12).Value = 7.505
$ws.Cells.Item(4, 13).Value = 7.991
$ws.Cells.Item(4, 14).Value = 2.177
$ws.Cells.Item(4, 15).Value = 6.767
$ws.Cells.Item(4, 16).Value = 9.513999999999999
$ws.Cells.Item(4, 17).Value = 5.969
$ws.Cells.Item(4, 18).Value = 0.9389999999999999
$ws.Cells.Item(4, 19).Value = 0.582
$ws.Cells.Item(4, 20).Value = 95.7
$ws.Cells.Item(4, 21).Value = 19.075
$ws.Cells.Item(4, 22).Value = 6.247
$ws.Cells.Item(4, 23).Value = 12.528
$ws.Cells.Item(4, 24).Value = 6.51
$ws.Cells.Item(4, 25).Value = 1.306
$ws.Cells.Item(4, 26).Value = 14.077
$ws.Cells.Item(4, 27).Value = 5.517
$ws.Cells.Item(4, 28).Value = 5.022
$ws.Cells.Item(4, 29).Value = 5.909
$ws.Cells.Item(4, 30).Value = 7.684
$ws.Cells.Item(4, 31).Value = 0.766
$ws.Cells.Item(4, 32).Value = 27.528
$ws.Cells.Item(4, 33).Value = 3.371
$ws.Cells.Item(4, 34).Value = 7.81
# Row 5
$ws.Cells.Item(5, 1).Value = 45099.52777777778
$ws.Cells.Item(5, 2).Value = 10.04
$ws.Cells.Item(5, 3).Value = 7.29
$ws.Cells.Item(5, 4).Value = 0.89
$ws.Cells.Item(5, 5).Value = 22.07
$ws.Cells.Item(5, 6).Value = 17.52
$ws.Cells.Item(5, 7).Value = 7.83
$ws.Cells.Item(5, 8).Value = 30.95
$ws.Cells.Item(5, 9).Value = 12.22
$ws.Cells.Item(5, 10).Value = 5.32
$ws.Cells.Item(5, 11).Value = 7.76
$ws.Cells.Item(5, 12).Value = 8.779999999999999
$ws.Cells.Item(5, 13).Value = 9.35
$ws.Cells.Item(5, 14).Value = 2.54
$ws.Cells.Item(5, 15).Value = 7.9
$ws.Cells.Item(5, 16).Value = 11.14
$ws.Cells.Item(5, 17).Value = 6.83
$ws.Cells.Item(5, 18).Value = 0.73
$ws.Cells.Item(5, 19).Value = 0.53
$ws.Cells.Item(5, 20).Value = 112.84
$ws.Cells.Item(5, 21).Value = 22.07
$ws.Cells.Item(5, 22).Value = 7.29
$ws.Cells.Item(5, 23).Value = 14.64
$ws.Cells.Item(5, 24).Value = 7.67
$ws.Cells.Item(5, 25).Value = 1.35
$ws.Cells.Item(5, 26).Value = 14.89
$ws.Cells.Item(5, 27).Value = 6.44
$ws.Cells.Item(5, 28).Value = 5.79
$ws.Cells.Item(5, 29).Value = 6.81
$ws.Cells.Item(5, 30).Value = 9.08
$ws.Cells.Item(5, 31).Value = 0.55
$ws.Cells.Item(5, 32).Value = 28.15
$ws.Cells.Item(5, 33).Value = 4.01
$ws.Cells.Item(5, 34).Value = 9.109999999999999
